$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.950543999671936
$ws.Range("B1").Value = 3.233506202697754
$ws.Range("C1").Value = 2.802427053451538
$ws.Range("D1").Value = 1.557803273200989
$ws.Range("E1").Value = 1.186467885971069
